$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix floating point precision on existing row 57 date value
$ws.Cells.Item(57, 1).Value = 44370.76653535417

# Add new row 58 data
$ws.Cells.Item(58, 1).Value = 44371.76033399082
$ws.Cells.Item(58, 1).NumberFormat = "yyyy-mm-dd HH:mm:ss UTC"
$ws.Cells.Item(58, 2).Value = 79009
$ws.Cells.Item(58, 3).Value = 66409
$ws.Cells.Item(58, 4).Value = 3541
$ws.Cells.Item(58, 5).Value = 2141
$ws.Cells.Item(58, 6).Value = 1534
$ws.Cells.Item(58, 7).Value = 20881
$ws.Cells.Item(58, 8).Value = 1453
$ws.Cells.Item(58, 9).Value = 893
$ws.Cells.Item(58, 10).Value = 182
